$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.812.52'
$ws.Range('E2').Value = '  +3.65%  '

$ws.Range('D3').Value = '2.424.25'
$ws.Range('E3').Value = '  +3.21%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.40'
$ws.Range('E5').Value = '  +1.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.34'
$ws.Range('E6').Value = '  +4.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +4.18%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.106'
$ws.Range('E9').Value = '  +1.39%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.71'
$ws.Range('E10').Value = '  +3.25%  '

$ws.Range('E11').Value = '  -1.87%  '

$ws.Range('E12').Value = '  +1.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.75'
$ws.Range('E13').Value = '  +4.24%  '

$ws.Range('D14').Value = '2.856.02'
$ws.Range('E14').Value = '  +3.16%  '

$ws.Range('D15').Value = '59.752.27'
$ws.Range('E15').Value = '  +3.64%  '

$ws.Range('E16').Value = '  +1.62%  '

$ws.Range('D17').Value = '2.424.18'
$ws.Range('E17').Value = '  +3.15%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.28'
$ws.Range('E18').Value = '  +2.83%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.38'
$ws.Range('E19').Value = '  +2.29%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '330.66'
$ws.Range('E20').Value = '  +0.80%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.69'
$ws.Range('E21').Value = '  -2.63%  '

$ws.Range('E22').Value = '  +0.24%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.75'
$ws.Range('E23').Value = '  +3.95%  '

$ws.Range('E24').Value = '  +3.33%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.61'
$ws.Range('E25').Value = '  +5.32%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('E27').Value = '  +3.14%  '

$ws.Range('D28').Value = '0.0₃0777'
$ws.Range('E28').Value = '  +5.74%  '

$ws.Range('E29').Value = '  +1.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.35'
$ws.Range('E30').Value = '  -0.40%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  +0.85%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.66'
$ws.Range('E32').Value = '  +1.93%  '

$ws.Range('E33').Value = '  +1.35%  '

$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.29'
$ws.Range('E35').Value = '  +5.26%  '

$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  +1.67%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.61'
$ws.Range('E38').Value = '  +1.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.40'
$ws.Range('E39').Value = '  +0.93%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '314.02'
$ws.Range('E40').Value = '  +8.90%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.408'
$ws.Range('E41').Value = '  -1.15%  '

$ws.Range('E42').Value = '  +0.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '138.88'
$ws.Range('E43').Value = '  -2.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0970'
$ws.Range('E44').Value = '  +2.57%  '

$ws.Range('E45').Value = '  +1.35%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.54'
$ws.Range('E46').Value = '  +5.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.580'
$ws.Range('E47').Value = '  +2.94%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0225'
$ws.Range('E48').Value = '  +2.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.391'
$ws.Range('E49').Value = '  +0.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.60'
$ws.Range('E50').Value = '  +1.34%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.06'
$ws.Range('E51').Value = '  -0.16%  '
